$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.791.66"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.343.14"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "514.67"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "133.85"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  -0.30%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.534"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.61%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "5.40"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.40%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.152"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.341"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.51%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "24.02"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "2.763.27"
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "56.807.66"
$ws.Range("E15").Value = "  +0.05%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000134"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "2.367.98"
$ws.Range("E17").Value = "  +1.96%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "10.44"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.66%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "326.71"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.50%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.19"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.17%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.69"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("E22").Value = "  +0.11%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "61.04"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.165"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "8.64"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +11.79%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.31"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +9.09%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "167.68"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0730"
$ws.Range("E29").Value = "  -0.06%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.68"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.19"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.31%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "18.45"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("E33").Value = "  -0.02%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.28"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.93%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.00"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.01%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.894"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.83%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.58"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.99%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "38.68"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.35%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "150.16"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +8.82%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.377"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.38%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.61"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.57%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "282.57"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.56%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.25"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.63%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0930"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0503"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.558"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.24%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "18.36"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +7.85%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0216"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.69%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "17.30"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.32%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "10.98"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.08%  "
